$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.727.27"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "2.379.02"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.639"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.631"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0923"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.39%  "
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "2.730.59"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "2.373.69"
$ws.Range("E17").Value = "  -2.82%  "
$ws.Range("D18").Value = "42.700.47"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000107"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +11.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "272.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.87%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "23.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  -2.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0910"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "35.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("E36").Value = "  +2.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0359"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.47%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.15%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.35%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.106"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.60%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.233"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.38%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.33%  "
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "118.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +31.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("D51").Value = "1.587.31"
$ws.Range("E51").Value = "  +6.31%  "
